# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-12-29 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-12-30 Saturday", 2)

# Update the worksheet answer cells by absolute (row, column) position so that
# values which coincide with other rows' old/new text cannot be mismatched by
# a text search.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "15÷8=1, 7"
$t.Cell(1, 2).Range.Text = "50÷5=10, 0"
$t.Cell(1, 3).Range.Text = "62÷4=15, 2"
$t.Cell(1, 4).Range.Text = "84÷3=28, 0"
$t.Cell(1, 5).Range.Text = "52÷7=7, 3"

$t.Cell(5, 1).Range.Text = "69÷6=11, 3"
$t.Cell(5, 2).Range.Text = "83÷7=11, 6"
$t.Cell(5, 3).Range.Text = "96÷6=16, 0"
$t.Cell(5, 4).Range.Text = "89÷2=44, 1"
$t.Cell(5, 5).Range.Text = "97÷4=24, 1"

$t.Cell(9, 1).Range.Text = "11÷3=3, 2"
$t.Cell(9, 2).Range.Text = "87÷5=17, 2"
$t.Cell(9, 3).Range.Text = "21÷6=3, 3"
$t.Cell(9, 4).Range.Text = "40÷5=8, 0"
$t.Cell(9, 5).Range.Text = "27÷9=3, 0"

$t.Cell(13, 1).Range.Text = "47÷4=11, 3"
$t.Cell(13, 2).Range.Text = "72÷9=8, 0"
$t.Cell(13, 3).Range.Text = "50÷7=7, 1"
$t.Cell(13, 4).Range.Text = "94÷4=23, 2"
$t.Cell(13, 5).Range.Text = "50÷3=16, 2"

$t.Cell(17, 1).Range.Text = "26÷5=5, 1"
$t.Cell(17, 2).Range.Text = "76÷9=8, 4"
$t.Cell(17, 3).Range.Text = "26÷7=3, 5"
$t.Cell(17, 4).Range.Text = "76÷9=8, 4"
$t.Cell(17, 5).Range.Text = "44÷7=6, 2"
